$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the DATE_OF_DISABLEMENT value for row 2 (still keeps its date style)
$ws.Range("L2").ClearContents()

# Update the DATE_OF_DISABLEMENT value for row 3 to 2017-10-01 (serial 43009)
$ws.Range("L3").Value = (Get-Date -Year 2017 -Month 10 -Day 1 -Hour 0 -Minute 0 -Second 0).Date

# Reflect the resulting selection/scroll position used when the file was saved
$ws.Activate()
$ws.Range("L3").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
